# Weekly update: insert two new price-report rows (a new "Jengibre" week)
# ahead of the existing row 28, pushing all subsequent rows down by two.
# This mirrors the diff: dimension grows from A1:R60 to A1:R62, and the
# data previously on rows 28-60 now lives on rows 30-62 unchanged, while
# the new rows 28 and 29 hold the newly reported week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 28; everything below (old rows 28-60)
# shifts down to rows 30-62, and column D keeps its date style because it
# is carried over from the row being pushed down.
$ws.Rows("28:29").Insert()

# New row 28
$ws.Cells.Item(28, 1).Value = 6
$ws.Cells.Item(28, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44571
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 100114007
$ws.Cells.Item(28, 7).Value = "Jengibre"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 140
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14571
$ws.Cells.Item(28, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 1121
$ws.Cells.Item(28, 17).Value = 13
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# New row 29
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 44571
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100114007
$ws.Cells.Item(29, 7).Value = "Jengibre"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 30
$ws.Cells.Item(29, 11).Value = 12000
$ws.Cells.Item(29, 12).Value = 12000
$ws.Cells.Item(29, 13).Value = 12000
$ws.Cells.Item(29, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 923
$ws.Cells.Item(29, 17).Value = 13
$ws.Cells.Item(29, 18).Value = "Hortaliza"
